# Update the "solidarity_support_incl_info_mean" data with the final
# prepared/rendered values (run prepare & render with final data).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.750480905905658
$ws.Range("K2").Value = 0.540833488586868
$ws.Range("L2").Value = 0.684994427052208
$ws.Range("N2").Value = 0.691916953645142

$ws.Range("B3").Value = 0.55247265659607
$ws.Range("K3").Value = 0.3844769302497
$ws.Range("L3").Value = 0.572679837389163
$ws.Range("N3").Value = 0.484051126996374

$ws.Range("B4").Value = 0.480079449338372
$ws.Range("K4").Value = 0.299554385931615
$ws.Range("L4").Value = 0.649188481861327
$ws.Range("N4").Value = 0.395528166509857

$ws.Range("B5").Value = 0.378764242866997
$ws.Range("K5").Value = 0.11169228057918
$ws.Range("L5").Value = 0.629942389379855
$ws.Range("M5").Value = 0.779759303305367
$ws.Range("N5").Value = 0.326278157947098

$ws.Range("B6").Value = 0.351543732037051
$ws.Range("K6").Value = 0.235572583007281
$ws.Range("L6").Value = -0.0291282837527572
$ws.Range("N6").Value = 0.283454808536633

$ws.Range("B7").Value = 0.35025634310265
$ws.Range("K7").Value = 0.0778112996303359
$ws.Range("L7").Value = 0.678885502636161
$ws.Range("M7").Value = 0.732044084797925
$ws.Range("N7").Value = 0.228454175346082

$ws.Range("B8").Value = 0.343504471736926
$ws.Range("K8").Value = 0.0841964845636908
$ws.Range("L8").Value = 0.350622027830543
$ws.Range("M8").Value = 0.627551511945665
$ws.Range("N8").Value = 0.301628017683413

$ws.Range("B9").Value = 0.341548248462238
$ws.Range("K9").Value = 0.196407827310263
$ws.Range("L9").Value = 0.380345366744009
$ws.Range("N9").Value = 0.301565761377154

$ws.Range("B10").Value = 0.314067189876321
$ws.Range("K10").Value = 0.340271717767617
$ws.Range("L10").Value = -0.023705423678396
$ws.Range("N10").Value = 0.241549262498416

$ws.Range("B11").Value = 0.0189523977290343
$ws.Range("K11").Value = -0.0721523686127008
$ws.Range("L11").Value = 0.0139318769794177
$ws.Range("N11").Value = -0.0376248011861876
